$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D16").Copy($ws.Range("B16"))
$ws.Range("B16").Value = "home.html"
$ws.Hyperlinks.Add($ws.Range("B16"), "home.html")
